$d = $word.ActiveDocument

# Locate the run of text that needs to be split: "one-year study of microhabitat
# and mesohabitat in Feather Rive" -> "two" + "-year study of microhabitat and
# mesohabitat in Feather Rive"
$target = $d.Content
$target.Find.Execute(
    "one-year study of microhabitat and mesohabitat in Feather Rive",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

$wholeStart = $target.Start
$wholeEnd = $target.End

# "one" occupies the first three characters of the match.
$oneStart = $wholeStart
$oneEnd = $wholeStart + 3

# Replace "one" with "two".
$rngOne = $d.Range($oneStart, $oneEnd)
$rngOne.Text = "two"

# "two" is also 3 characters long, so the remainder ("-year study of
# microhabitat and mesohabitat in Feather Rive") now starts right after it
# and ends where the original match ended.
$tailStart = $oneStart + 3
$tailEnd = $wholeEnd

# Re-stamp formatting on the tail so it becomes its own run, separate from
# the "two" run that precedes it.
$rngTail = $d.Range($tailStart, $tailEnd)
$rngTail.Font.Name = "Calibri"
$rngTail.Font.NameBi = "Calibri"

# Re-stamp formatting on "two" last, so it splits away from the preceding
# (unchanged) run as well, leaving it as its own standalone run.
$rngTwo = $d.Range($oneStart, $tailStart)
$rngTwo.Font.Name = "Calibri"
$rngTwo.Font.NameBi = "Calibri"
